$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Hgf"
$ws.Range("C2").Value = "Sdc1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1088396666666667
$ws.Range("H2").Value = 0.326519
$ws.Range("I2").Value = 0.002750770615347974
$ws.Range("J2").Value = 0.002750770615347974
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.378475333333333
$ws.Range("N2").Value = 4.135426
$ws.Range("O2").Value = 0.05609715574531157
$ws.Range("P2").Value = 0.05609715574531156
$ws.Range("Q2").Value = 0.1500327957882222
$ws.Range("R2").Value = 1.350295162094
$ws.Range("S2").Value = 0.0001543104076288018
$ws.Range("T2").Value = 0.0001543104076288018

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Hgf"
$ws.Range("C3").Value = "Sdc1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1088396666666667
$ws.Range("H3").Value = 0.326519
$ws.Range("I3").Value = 0.002750770615347974
$ws.Range("J3").Value = 0.002750770615347974
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 13.16176133333333
$ws.Range("N3").Value = 39.485284
$ws.Range("O3").Value = 0.5356188518899525
$ws.Range("P3").Value = 0.5356188518899525
$ws.Range("Q3").Value = 1.432521716266222
$ws.Range("R3").Value = 12.892695446396
$ws.Range("S3").Value = 0.0014733645988053
$ws.Range("T3").Value = 0.0014733645988053

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Hgf"
$ws.Range("C4").Value = "Sdc1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1088396666666667
$ws.Range("H4").Value = 0.326519
$ws.Range("I4").Value = 0.002750770615347974
$ws.Range("J4").Value = 0.002750770615347974
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.829094
$ws.Range("N4").Value = 2.487282
$ws.Range("O4").Value = 0.03374004171190829
$ws.Range("P4").Value = 0.03374004171190828
$ws.Range("Q4").Value = 0.09023831459533334
$ws.Range("R4").Value = 0.812144831358
$ws.Range("S4").Value = 0.00009281111530173226
$ws.Range("T4").Value = 0.00009281111530173225

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Hgf"
$ws.Range("C5").Value = "Sdc1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1088396666666667
$ws.Range("H5").Value = 0.326519
$ws.Range("I5").Value = 0.002750770615347974
$ws.Range("J5").Value = 0.002750770615347974
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.203668
$ws.Range("N5").Value = 27.611004
$ws.Range("O5").Value = 0.3745439506528278
$ws.Range("P5").Value = 0.3745439506528276
$ws.Range("Q5").Value = 1.001724157230667
$ws.Range("R5").Value = 9.015517415076001
$ws.Range("S5").Value = 0.00103028449361214
$ws.Range("T5").Value = 0.00103028449361214

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Hgf"
$ws.Range("C6").Value = "Sdc1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 15.07419333333333
$ws.Range("H6").Value = 45.22258
$ws.Range("I6").Value = 0.3809791902285103
$ws.Range("J6").Value = 0.3809791902285103
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.378475333333333
$ws.Range("N6").Value = 4.135426
$ws.Range("O6").Value = 0.05609715574531157
$ws.Range("P6").Value = 0.05609715574531156
$ws.Range("Q6").Value = 20.77940367989778
$ws.Range("R6").Value = 187.01463311908
$ws.Range("S6").Value = 0.02137184896997143
$ws.Range("T6").Value = 0.02137184896997142

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Hgf"
$ws.Range("C7").Value = "Sdc1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 15.07419333333333
$ws.Range("H7").Value = 45.22258
$ws.Range("I7").Value = 0.3809791902285103
$ws.Range("J7").Value = 0.3809791902285103
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.16176133333333
$ws.Range("N7").Value = 39.485284
$ws.Range("O7").Value = 0.5356188518899525
$ws.Range("P7").Value = 0.5356188518899525
$ws.Range("Q7").Value = 198.4029349458578
$ws.Range("R7").Value = 1785.62641451272
$ws.Range("S7").Value = 0.2040596364641585
$ws.Range("T7").Value = 0.2040596364641585

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Hgf"
$ws.Range("C8").Value = "Sdc1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 15.07419333333333
$ws.Range("H8").Value = 45.22258
$ws.Range("I8").Value = 0.3809791902285103
$ws.Range("J8").Value = 0.3809791902285103
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.829094
$ws.Range("N8").Value = 2.487282
$ws.Range("O8").Value = 0.03374004171190829
$ws.Range("P8").Value = 0.03374004171190828
$ws.Range("Q8").Value = 12.49792324750667
$ws.Range("R8").Value = 112.48130922756
$ws.Range("S8").Value = 0.01285425376967898
$ws.Range("T8").Value = 0.01285425376967898

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Hgf"
$ws.Range("C9").Value = "Sdc1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 15.07419333333333
$ws.Range("H9").Value = 45.22258
$ws.Range("I9").Value = 0.3809791902285103
$ws.Range("J9").Value = 0.3809791902285103
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.203668
$ws.Range("N9").Value = 27.611004
$ws.Range("O9").Value = 0.3745439506528278
$ws.Range("P9").Value = 0.3745439506528276
$ws.Range("Q9").Value = 138.7378708078134
$ws.Range("R9").Value = 1248.64083727032
$ws.Range("S9").Value = 0.1426934510247015
$ws.Range("T9").Value = 0.1426934510247014

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Hgf"
$ws.Range("C10").Value = "Sdc1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 23.69325166666667
$ws.Range("H10").Value = 71.079755
$ws.Range("I10").Value = 0.5988138558556568
$ws.Range("J10").Value = 0.5988138558556569
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.378475333333333
$ws.Range("N10").Value = 4.135426
$ws.Range("O10").Value = 0.05609715574531157
$ws.Range("P10").Value = 0.05609715574531156
$ws.Range("Q10").Value = 32.66056298895889
$ws.Range("R10").Value = 293.94506690063
$ws.Range("S10").Value = 0.03359175413438533
$ws.Range("T10").Value = 0.03359175413438533

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Hgf"
$ws.Range("C11").Value = "Sdc1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 23.69325166666667
$ws.Range("H11").Value = 71.079755
$ws.Range("I11").Value = 0.5988138558556568
$ws.Range("J11").Value = 0.5988138558556569
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 13.16176133333333
$ws.Range("N11").Value = 39.485284
$ws.Range("O11").Value = 0.5356188518899525
$ws.Range("P11").Value = 0.5356188518899525
$ws.Range("Q11").Value = 311.8449236472689
$ws.Range("R11").Value = 2806.60431282542
$ws.Range("S11").Value = 0.3207359899692024
$ws.Range("T11").Value = 0.3207359899692025

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Hgf"
$ws.Range("C12").Value = "Sdc1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 23.69325166666667
$ws.Range("H12").Value = 71.079755
$ws.Range("I12").Value = 0.5988138558556568
$ws.Range("J12").Value = 0.5988138558556569
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.829094
$ws.Range("N12").Value = 2.487282
$ws.Range("O12").Value = 0.03374004171190829
$ws.Range("P12").Value = 0.03374004171190828
$ws.Range("Q12").Value = 19.64393279732333
$ws.Range("R12").Value = 176.79539517591
$ws.Range("S12").Value = 0.0202040044742385
$ws.Range("T12").Value = 0.0202040044742385

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Hgf"
$ws.Range("C13").Value = "Sdc1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 23.69325166666667
$ws.Range("H13").Value = 71.079755
$ws.Range("I13").Value = 0.5988138558556568
$ws.Range("J13").Value = 0.5988138558556569
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 9.203668
$ws.Range("N13").Value = 27.611004
$ws.Range("O13").Value = 0.3745439506528278
$ws.Range("P13").Value = 0.3745439506528276
$ws.Range("Q13").Value = 218.0648221804467
$ws.Range("R13").Value = 1962.58339962402
$ws.Range("S13").Value = 0.2242821072778306
$ws.Range("T13").Value = 0.2242821072778306

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Hgf"
$ws.Range("C14").Value = "Sdc1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.6906883333333335
$ws.Range("H14").Value = 2.072065
$ws.Range("I14").Value = 0.01745618330048481
$ws.Range("J14").Value = 0.01745618330048481
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.378475333333333
$ws.Range("N14").Value = 4.135426
$ws.Range("O14").Value = 0.05609715574531157
$ws.Range("P14").Value = 0.05609715574531156
$ws.Range("Q14").Value = 0.9520968305211113
$ws.Range("R14").Value = 8.56887147469
$ws.Range("S14").Value = 0.0009792422333260032
$ws.Range("T14").Value = 0.0009792422333260032

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Hgf"
$ws.Range("C15").Value = "Sdc1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.6906883333333335
$ws.Range("H15").Value = 2.072065
$ws.Range("I15").Value = 0.01745618330048481
$ws.Range("J15").Value = 0.01745618330048481
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 13.16176133333333
$ws.Range("N15").Value = 39.485284
$ws.Range("O15").Value = 0.5356188518899525
$ws.Range("P15").Value = 0.5356188518899525
$ws.Range("Q15").Value = 9.090674999051112
$ws.Range("R15").Value = 81.81607499146001
$ws.Range("S15").Value = 0.009349860857786236
$ws.Range("T15").Value = 0.009349860857786236

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Hgf"
$ws.Range("C16").Value = "Sdc1"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.6906883333333335
$ws.Range("H16").Value = 2.072065
$ws.Range("I16").Value = 0.01745618330048481
$ws.Range("J16").Value = 0.01745618330048481
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.829094
$ws.Range("N16").Value = 2.487282
$ws.Range("O16").Value = 0.03374004171190829
$ws.Range("P16").Value = 0.03374004171190828
$ws.Range("Q16").Value = 0.5726455530366668
$ws.Range("R16").Value = 5.153809977330001
$ws.Range("S16").Value = 0.0005889723526890744
$ws.Range("T16").Value = 0.0005889723526890742

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Hgf"
$ws.Range("C17").Value = "Sdc1"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.6906883333333335
$ws.Range("H17").Value = 2.072065
$ws.Range("I17").Value = 0.01745618330048481
$ws.Range("J17").Value = 0.01745618330048481
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 9.203668
$ws.Range("N17").Value = 27.611004
$ws.Range("O17").Value = 0.3745439506528278
$ws.Range("P17").Value = 0.3745439506528276
$ws.Range("Q17").Value = 6.356866111473335
$ws.Range("R17").Value = 57.21179500326001
$ws.Range("S17").Value = 0.006538107856683499
$ws.Range("T17").Value = 0.006538107856683497

